# Updates the cryptos price/volume table (Sheet1) with refreshed
# coinranking.com figures. Price strings that look numeric ("0.999",
# "554.90", etc.) are written with NumberFormat="@" first and the
# style reset to "Normal" afterward so they stay plain text (matching
# the original inlineStr cells) instead of being auto-coerced into
# numbers and losing their formatting (trailing zeros, thousands dots).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.346.70'
$ws.Range("E2").Value = '  -5.77%  '
$ws.Range("D3").Value = '3.452.07'
$ws.Range("E3").Value = '  -7.46%  '
$r = $ws.Range("D4")
$r.NumberFormat = "@"
$r.Value = '0.999'
$r.Style = "Normal"
$ws.Range("E4").Value = '  -0.07%  '
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = '554.90'
$r.Style = "Normal"
$ws.Range("E5").Value = '  -9.30%  '
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = '179.60'
$r.Style = "Normal"
$ws.Range("E6").Value = '  -6.63%  '
$ws.Range("D7").Value = '3.448.61'
$ws.Range("E7").Value = '  -7.40%  '
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = '0.594'
$r.Style = "Normal"
$ws.Range("E8").Value = '  -6.90%  '
$ws.Range("E9").Value = '  +0.18%  '
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = '0.638'
$r.Style = "Normal"
$ws.Range("E10").Value = '  -12.58%  '
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = '0.139'
$r.Style = "Normal"
$ws.Range("E11").Value = '  -14.20%  '
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = '50.74'
$r.Style = "Normal"
$ws.Range("E12").Value = '  -16.42%  '
$ws.Range("E13").Value = '  -14.94%  '
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = '9.35'
$r.Style = "Normal"
$ws.Range("E14").Value = '  -12.62%  '
$ws.Range("D15").Value = '3.998.57'
$ws.Range("E15").Value = '  -7.46%  '
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = '0.124'
$r.Style = "Normal"
$ws.Range("E16").Value = '  -1.87%  '
$ws.Range("D17").Value = '3.451.65'
$ws.Range("E17").Value = '  -7.28%  '
$ws.Range("D18").Value = '64.977.48'
$ws.Range("E18").Value = '  -6.05%  '
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = '17.54'
$r.Style = "Normal"
$ws.Range("E19").Value = '  -10.11%  '
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = '11.53'
$r.Style = "Normal"
$ws.Range("E20").Value = '  -11.19%  '
$ws.Range("E21").Value = '  -11.89%  '
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = '373.11'
$r.Style = "Normal"
$ws.Range("E22").Value = '  -9.81%  '
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = '4.04'
$r.Style = "Normal"
$ws.Range("E23").Value = '  -11.80%  '
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = '81.90'
$r.Style = "Normal"
$ws.Range("E24").Value = '  -8.90%  '
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = '10.61'
$r.Style = "Normal"
$ws.Range("E25").Value = '  -3.00%  '
$ws.Range("E26").Value = '  -1.22%  '
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = '2.75'
$r.Style = "Normal"
$ws.Range("E27").Value = '  -10.33%  '
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = '11.69'
$r.Style = "Normal"
$ws.Range("E28").Value = '  -9.50%  '
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = '3.38'
$r.Style = "Normal"
$ws.Range("E29").Value = '  -11.51%  '
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = '8.44'
$r.Style = "Normal"
$ws.Range("E30").Value = '  -13.29%  '
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = '29.84'
$r.Style = "Normal"
$ws.Range("E31").Value = '  -10.13%  '
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = '7.03'
$r.Style = "Normal"
$ws.Range("E32").Value = '  -9.77%  '
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = '608.51'
$r.Style = "Normal"
$ws.Range("E33").Value = '  -4.47%  '
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = '11.70'
$r.Style = "Normal"
$ws.Range("E34").Value = '  -9.03%  '
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = '62.38'
$r.Style = "Normal"
$ws.Range("E35").Value = '  -5.67%  '
$ws.Range("E36").Value = '  -13.11%  '
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = '39.99'
$r.Style = "Normal"
$ws.Range("E37").Value = '  -13.04%  '
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = '1.00'
$r.Style = "Normal"
$ws.Range("E38").Value = '  +0.11%  '
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = '0.391'
$r.Style = "Normal"
$ws.Range("E39").Value = '  -6.71%  '
$ws.Range("E40").Value = '  -0.31%  '
$ws.Range("D41").Value = '0.0₃0704'
$ws.Range("E41").Value = '  -15.34%  '
$ws.Range("E42").Value = '  -9.51%  '
$ws.Range("D43").Value = '2.884.45'
$ws.Range("E43").Value = '  -0.23%  '
$ws.Range("E44").Value = '  -12.77%  '
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = '2.42'
$r.Style = "Normal"
$ws.Range("E45").Value = '  -8.51%  '
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = '3.09'
$r.Style = "Normal"
$ws.Range("E46").Value = '  -0.51%  '
$ws.Range("E47").Value = '  -13.63%  '
$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = '137.50'
$r.Style = "Normal"
$ws.Range("E48").Value = '  -4.52%  '
$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = '0.125'
$r.Style = "Normal"
$ws.Range("E49").Value = '  -10.80%  '
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = '2.43'
$r.Style = "Normal"
$ws.Range("E50").Value = '  -11.62%  '
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = '8.03'
$r.Style = "Normal"
$ws.Range("E51").Value = '  -12.46%  '
